# Update Lorenz distribution values on the "Sim_5yr" and "Sim_10yr" sheets
# Row 2 corresponds to age_bin "25-30" (and age_bin_5yr for sheet2), columns
# B:E are lorenz_20, lorenz_40, lorenz_60, lorenz_80 respectively.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Sim_5yr", "Sim_10yr")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("B2").Value = 0.004854158060199548
    $ws.Range("C2").Value = 0.02385602883565608
    $ws.Range("D2").Value = 0.06645552261722731
    $ws.Range("E2").Value = 0.1618653300067957
}
